{"js": "// \u267f\ufe0fA11Y: added <RadioSelector> forced colors behavior\n//\n// The underlying visible text of the first paragraph and of the\n// \"Determines the contrast...\" paragraph does not change \u2014 only the\n// run/proofErr structure around them is simplified (spell/grammar-check\n// markers removed, a word re-split across two runs). A brand-new\n// paragraph describing the \"forced colors\" override is appended right\n// after the contrast paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 1: \"This is version 10 of my portfolio website. ...\" ---\n// Collapse the proofErr-delimited runs (Vite / \"open source\") into a\n// single contiguous run with identical text.\nconst introPara = paragraphs.items[0];\nconst introText =\n  \"This is version 10 of my portfolio website. It is a single-page \" +\n  \"application built with React, bundled using Vite, and animated with \" +\n  \"Framer Motion. You can view the open source code on GitHub.\";\nintroPara.getRange().insertText(introText, \"Replace\");\n\n// --- Last paragraph: \"Determines the contrast between foreground ...\" ---\n// Re-split \"satisfies\" into \"satisf\" + \"ies\" and merge the remaining\n// trailing runs into one, using insertOoxml so the exact run\n// boundaries from the diff are produced (text content is unchanged).\nconst contrastPara = paragraphs.items[paragraphs.items.length - 1];\nconst contrastPackage =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  \"<w:r><w:t>Determines the contrast between foreground and background colors. \\u2018Default\\u2019 satisfies WCAG level AA and \\u2018high\\u2019 satisf</w:t></w:r>\" +\n  \"<w:r><w:t>ies</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> level AAA. \\u2018Auto\\u2019 will respect the system contrast.</w:t></w:r>' +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\ncontrastPara.insertOoxml(contrastPackage, \"Replace\");\n\n// --- New paragraph after the contrast paragraph ---\n// \"This setting is currently overridden as forced colors are active.\"\n// split across five runs, as produced by incremental typing/edits.\nconst newPara = contrastPara.insertParagraph(\"\", \"After\");\nconst newParaPackage =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">This setting is </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">currently </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">overridden </w:t></w:r>' +\n  \"<w:r><w:t>as</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> forced colors are active.</w:t></w:r>' +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\nnewPara.insertOoxml(newParaPackage, \"Replace\");\n\nawait context.sync();\n", "ps1": "# \u267f\ufe0fA11Y: added <RadioSelector> forced colors behavior\n#\n# The underlying visible text of the first paragraph and of the\n# \"Determines the contrast...\" paragraph does not change - only the\n# run/proofErr structure around them is simplified (spell/grammar-check\n# markers removed, a word re-split across two runs). A brand-new\n# paragraph describing the \"forced colors\" override is appended right\n# after the contrast paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: \"This is version 10 of my portfolio website. ...\" ---\n# Collapse the proofErr-delimited runs (Vite / \"open source\") into a\n# single contiguous run with identical text. Keep the paragraph's own\n# mark intact (only clear its contents) so the following empty\n# paragraph is not swallowed along with it.\n$introPara = $d.Paragraphs.Item(1)\n$introRange = $introPara.Range\n$introRange.MoveEnd(1, -1) | Out-Null\n$introRange.Delete()\n$introRange.InsertAfter(\"This is version 10 of my portfolio website. It is a single-page application built with React, bundled using Vite, and animated with Framer Motion. You can view the open source code on GitHub.\")\n\n# --- Last paragraph: \"Determines the contrast between foreground ...\" ---\n# Re-split \"satisfies\" into \"satisf\" + \"ies\" and merge the remaining\n# trailing runs into one. Delete the WHOLE paragraph (including its\n# mark - it is the last paragraph in the body, right before sectPr, so\n# nothing else is absorbed) and re-insert it via InsertXML so the exact\n# run boundaries from the diff are produced (visible text is unchanged).\n$contrastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$contrastRange = $contrastPara.Range\n$contrastRange.Delete()\n$contrastXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Determines the contrast between foreground and background colors. \u2018Default\u2019 satisfies WCAG level AA and \u2018high\u2019 satisf</w:t></w:r><w:r><w:t>ies</w:t></w:r><w:r><w:t xml:space=\"preserve\"> level AAA. \u2018Auto\u2019 will respect the system contrast.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$contrastRange.InsertXML($contrastXml)\n\n# --- New paragraph after the contrast paragraph ---\n# \"This setting is currently overridden as forced colors are active.\"\n# split across five runs, as produced by incremental typing/edits.\n$contrastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$contrastPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newRange = $newPara.Range\n$newRange.Delete()\n$newXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">This setting is </w:t></w:r><w:r><w:t xml:space=\"preserve\">currently </w:t></w:r><w:r><w:t xml:space=\"preserve\">overridden </w:t></w:r><w:r><w:t>as</w:t></w:r><w:r><w:t xml:space=\"preserve\"> forced colors are active.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$newRange.InsertXML($newXml)\n"}
